$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.132.81"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.907.28"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "328.02"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4631"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "0.3933"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "46.86"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "0.07984"
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("D11").Value = "1.006"
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("D12").Value = "22.38"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "2.059.85"
$ws.Range("E13").Value = "  +11.76%  "
$ws.Range("D14").Value = "7.155"
$ws.Range("E14").Value = "  +2.56%  "
$ws.Range("D15").Value = "5.774"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "0.06987"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "88.93"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "0.00001012"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "17.27"
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "29.208.63"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("D23").Value = "5.393"
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "2.296.74"
$ws.Range("E25").Value = "  +10.94%  "
$ws.Range("E26").Value = "  -3.23%  "
$ws.Range("D27").Value = "156.64"
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("D29").Value = "5.896"
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("D30").Value = "2.009"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "119.86"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "0.09363"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "0.9294"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "5.370"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("D35").Value = "1.344"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "3.279"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("D37").Value = "0.05845"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "8.057"
$ws.Range("E38").Value = "  +3.74%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.159"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("D40").Value = "0.02102"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5780"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.1814"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "10.00"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "2.269"
$ws.Range("E44").Value = "  +11.26%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "11.98"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5431"
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.07138"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.889"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "1.119"
$ws.Range("E49").Value = "  -3.95%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "112.42"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").Value = "2.518"
$ws.Range("E51").Value = "  +6.47%  "
